# Adds two new summary columns to the right of the existing data:
#   GN = "median" of the party-difference values on each row
#   GO = "mean"   of the party-difference values on each row
# (These are per-row cross-topic statistics computed upstream; the values
#  below are the authoritative numbers produced by that calculation.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("GN1").Value = "median"
$ws.Range("GO1").Value = "mean"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered) by copying the format from the adjacent header cell GM1.
$ws.Range("GM1").Copy()
$ws.Range("GN1:GO1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (years 1972-2024, sheet rows 2-39) -----------------------
$medians = @(0.03285440856456028, 0.02755272203484037, 0.01965076489533008, 0.03408115413539169, 0.02785415740270621, 0.01937860433113033, 0.0254277931161527, 0.03633143445512643, 0.04532672586015944, 0.03274086227326588, 0.04267665892590883, 0.04349757932343196, 0.0393070092860975, 0.03067360192769471, 0.02536763926865549, 0.03662211388404191, 0.04098778455393001, 0.04796064267122746, 0.0283096857746935, 0.04479098901783383, 0.0345048704127461, 0.04328053978814657, 0.03101679295986759, 0.02983280177363729, 0.04155238736148028, 0.03794947492770474, 0.04176630428601012, 0.04396626545484317, 0.05947401499433308, 0.05173806730446739, 0.04785520980266225, 0.05979293514149104, 0.069765323146046, 0.06883674994238953, 0.07678965979778807, 0.0702439680357677, 0.07031580277258753, 0.07203577175959819)

$means = @(0.02853783047579867, 0.0401034855044566, 0.02418650682259138, 0.04216732468472677, 0.03314734242471936, 0.03066725527911303, 0.03277633436352156, 0.03828066581938715, 0.05264547661864068, 0.03815921276152804, 0.04734141505915695, 0.04546515254080564, 0.05201809402526514, 0.0469421887066003, 0.038955425437353, 0.04133935008587139, 0.042771459998787, 0.04862403392917938, 0.04004979720575044, 0.0553213337299854, 0.04138113423598297, 0.04713347507228979, 0.0419735020549061, 0.03871561571854708, 0.05207445293199564, 0.04771070083594445, 0.05220129192850576, 0.05529598297025578, 0.06047632659264623, 0.05798589891502316, 0.06068782307594302, 0.07234216568031356, 0.08131974350379645, 0.08745639917045274, 0.09467446729125548, 0.09974756583769276, 0.1026572806889232, 0.1089851809530814)

for ($i = 0; $i -lt $medians.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 196).Value = $medians[$i]  # column GN
    $ws.Cells.Item($row, 197).Value = $means[$i]    # column GO
}
